$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ((Intercept))
$ws.Range("B2").Value = 142810.25641
$ws.Range("D2").Value = 138.179495

# Row 3 (household_group_collapsed)
$ws.Range("B3").Value = 18084.358538
$ws.Range("D3").Value = 8.748977999999999
$ws.Range("E3").Value = 0.00022

# Row 4 (Residuals)
$ws.Range("B4").Value = 230473.321037
$ws.Range("C4").Value = 223

# Row 5 (SM-Control)
$ws.Range("G5").Value = -8.094238000000001
$ws.Range("H5").Value = -22.134558
$ws.Range("I5").Value = 5.946082
$ws.Range("J5").Value = 0.36364

# Row 6 (SM + Traps-Control)
$ws.Range("G6").Value = 12.166331
$ws.Range("H6").Value = -2.952121
$ws.Range("I6").Value = 27.284783
$ws.Range("J6").Value = 0.141411

# Row 7 (SM + Traps-SM)
$ws.Range("G7").Value = 20.260569
$ws.Range("H7").Value = 8.830862
$ws.Range("I7").Value = 31.690276
$ws.Range("J7").Value = 0.000122
